# Apply the "new .ttl from Google sheet has been generated" update.
#
# Summary of the change (derived from the OOXML diff):
#  - Row 13 (old "dct:creator" / "Jitka") is removed; everything below it
#    shifts up by one row, which also removes the trailing blank row 90
#    (dimension goes from A1:T90 to A1:T89).
#  - A handful of cell values are updated to reflect the newly generated
#    vocabulary metadata and the new "vars:working" term row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 13 entirely; this shifts rows 14-90 up by one (matching the
# diff, since rows 14-90 in the original file are identical, cell for
# cell, to rows 13-89 after the edit) and drops the now-superfluous last
# blank row, updating the sheet dimension from A1:T90 to A1:T89.
$ws.Rows(13).Delete()

# ConceptScheme PREFIX URI for "vars"
$ws.Range("C3").Value = "http://ontology.deic.dk/cv/DTUbib-vocab/"

# Vocabulary metadata block
$ws.Range("B10").Value = "DTUbib test"
$ws.Range("B11").Value = "DTUbib test"
$ws.Range("B12").Value = "Hannah"

# New term rows at the bottom of the term definitions
$ws.Range("A19").Value = "vars:test"
$ws.Range("B19").Value = "test"

$ws.Range("A20").Value = "vars:working"
$ws.Range("B20").Value = "working"
